{"js": "// Revise the \"elevator pitch\" paragraph: rewrite the intro/topics sentence,\n// drop \"for improvisation\" before \"and working on music\", and replace the\n// closing sentence about movement/gesture with the embodiment/mapping one.\n// Applied as a sequence of unambiguous, ordered find-and-replace passes so\n// the surrounding (unchanged) text is left untouched.\n\nconst replacements = [\n  {\n    find:\n      \"My research will focus on interactive computer music and gestural controllers. I intend to explore the historical/cultural background, and diverse or idiosyncratic approaches to mapping, as well as touch a bit on the phenomenology and aesthetics\",\n    replace:\n      \"I will continue researching about interactive computer music and gestural controllers. Covering topics such as early historical developments, idiosyncratic approaches to mapping, and a bit on the phenomenology\",\n  },\n  {\n    find: \" of human-computer interaction. My practice will focus on \",\n    replace:\n      \" of human-computer interaction. Machine learning will be explored as a generative technique for both mapping and cross-synthesis. I will relate this research to my own practice, which will focus on \",\n  },\n  {\n    find:\n      \"performance of interactive systems for improvisation and working on music\",\n    replace: \"performance of interactive systems and working on music\",\n  },\n  {\n    find:\n      \"exploring sonically the way movement and gesture leaves physical traces and influence perception, and exploring constraints imposed by situatedness. \",\n    replace:\n      \"exploring the way embodiment is influenced by the controller/mapping combination, working within self-imposed constraints, and employing meaningful categories of gesture.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Revise the \"elevator pitch\" paragraph: rewrite the intro/topics sentence,\n# drop \"for improvisation\" before \"and working on music\", and replace the\n# closing sentence about movement/gesture with the embodiment/mapping one.\n# Applied as a sequence of unambiguous, ordered Find/Replace passes so the\n# surrounding (unchanged) text is left untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find    = \"My research will focus on interactive computer music and gestural controllers. I intend to explore the historical/cultural background, and diverse or idiosyncratic approaches to mapping, as well as touch a bit on the phenomenology and aesthetics\"\n        Replace = \"I will continue researching about interactive computer music and gestural controllers. Covering topics such as early historical developments, idiosyncratic approaches to mapping, and a bit on the phenomenology\"\n    },\n    @{\n        Find    = \" of human-computer interaction. My practice will focus on \"\n        Replace = \" of human-computer interaction. Machine learning will be explored as a generative technique for both mapping and cross-synthesis. I will relate this research to my own practice, which will focus on \"\n    },\n    @{\n        Find    = \"performance of interactive systems for improvisation and working on music\"\n        Replace = \"performance of interactive systems and working on music\"\n    },\n    @{\n        Find    = \"exploring sonically the way movement and gesture leaves physical traces and influence perception, and exploring constraints imposed by situatedness. \"\n        Replace = \"exploring the way embodiment is influenced by the controller/mapping combination, working within self-imposed constraints, and employing meaningful categories of gesture.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, 0, $false, $r.Replace, 2) | Out-Null\n}\n"}
